# The deck currently uses the "Integral" theme (ppt/theme/theme2.xml,
# wired to the slide master / presentation) while the notes master is
# wired to ppt/theme/theme1.xml, which still holds the stock "Office
# Theme" color scheme. The authored change swaps the two themes so the
# slide master (and therefore every slide) renders with the Office
# Theme palette instead of Integral.
#
# PowerPoint's COM object model only exposes theme colors through
# ColorScheme.Colors(index).RGB (there's no supported way to swap which
# theme part a master points at, or to rename a theme). Each index maps
# to one of the 12 <a:clrScheme> children, in document order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2
#   7 accent3  8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
# RGB is stored as a PowerPoint "long" (0x00BBGGRR), so convert each
# target hex color before assigning it.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

$cs.Colors(1).RGB  = RGB(0x00, 0x00, 0x00)   # dk1      -> 000000
$cs.Colors(2).RGB  = RGB(0xFF, 0xFF, 0xFF)   # lt1      -> FFFFFF
$cs.Colors(3).RGB  = RGB(0x44, 0x54, 0x6A)   # dk2      -> 44546A
$cs.Colors(4).RGB  = RGB(0xE7, 0xE6, 0xE6)   # lt2      -> E7E6E6
$cs.Colors(5).RGB  = RGB(0x5B, 0x9B, 0xD5)   # accent1  -> 5B9BD5
$cs.Colors(6).RGB  = RGB(0xED, 0x7D, 0x31)   # accent2  -> ED7D31
$cs.Colors(7).RGB  = RGB(0xA5, 0xA5, 0xA5)   # accent3  -> A5A5A5
$cs.Colors(8).RGB  = RGB(0xFF, 0xC0, 0x00)   # accent4  -> FFC000
$cs.Colors(9).RGB  = RGB(0x44, 0x72, 0xC4)   # accent5  -> 4472C4
$cs.Colors(10).RGB = RGB(0x70, 0xAD, 0x47)   # accent6  -> 70AD47
$cs.Colors(11).RGB = RGB(0x05, 0x63, 0xC1)   # hlink    -> 0563C1
$cs.Colors(12).RGB = RGB(0x95, 0x4F, 0x72)   # folHlink -> 954F72
